# Resident Services_Requirements.xlsx - "Add files via upload" edit
#
# Adds a new "Reg Proc" column (the 20th column, header cell T2) to the
# "Table2" Excel table on the "Details" sheet, and fills in clarification /
# research notes for several rows in that new column (plus an update to the
# existing "Comments" (S) cell for row 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")
$lo = $ws.ListObjects.Item(1)

# --- 1. Grow the table by one column and name it ------------------------
$newCol = $lo.ListColumns.Add()
$hdrCell = $lo.HeaderRowRange.Cells.Item(1, 20)
$hdrCell.Value = "Reg Proc"

# --- 2. Fill in the new column's data cells (and the updated S8 comment) -
# Values are written in the same order the source workbook first introduced
# them in, so freshly-created shared-string entries line up with the target.
$ws.Range("T5").Value = "No Mapping of such kind from Reg Processor`nID Repo- Might not be there in ID Repo as well"
$ws.Range("T6").Value = "ID Repo- need to know "
$ws.Range("T9").Value = "Under processing`nProcessed"
$ws.Range("T8").Value = "Under processing`nProcessed`n"
$ws.Range("S8").Value = "Reg proc`nArchival policy"
$ws.Range("T10").Value = "E-UIN Generation"
$ws.Range("T7").Value = "there shud be a label as Res_Service`nReg Client packet needs to be understood`nService from Reg proc needs to be developed"
$ws.Range("T4").Value = "When UIN IS needed to be generated`n1.the Acknowledgment from Print queue- what needs to be done`nTime period `n2. If there is a print failure- no need to handle from MOSIP`nUser Story ?"

# --- 3. Match formatting of the rest of the table ------------------------
# Header cell: copy the look of the neighbouring header cell (centered,
# italic, bordered) instead of the plain default style.
$ws.Range("S2").Copy()
$ws.Range("T2").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Wrapped, left/top aligned cells (matches the other long free-text cells).
foreach ($addr in @("T4", "T5", "T7", "S8", "T8", "T9")) {
    $c = $ws.Range($addr)
    $c.HorizontalAlignment = -4131   # xlLeft
    $c.VerticalAlignment = -4160     # xlTop
    $c.WrapText = $true
}

# Plain, left/top aligned, non-wrapping cells (short single-line values).
foreach ($addr in @("T6", "T10")) {
    $c = $ws.Range($addr)
    $c.HorizontalAlignment = -4131   # xlLeft
    $c.VerticalAlignment = -4160     # xlTop
    $c.WrapText = $false
}

# --- 4. Column width for the new column -----------------------------------
$ws.Columns.Item(20).ColumnWidth = 31.2

# --- 5. Update the active selection to the newly edited cell --------------
$ws.Activate()
$ws.Range("T4").Select()
